{"js": "// Update the worksheet date and the 25 division problems to the new values.\nconst replacements = [\n  [\"2026-01-27 Tuesday\", \"2026-01-28 Wednesday\"],\n  [\"441\u00f72=\", \"143\u00f78=\"],\n  [\"145\u00f72=\", \"808\u00f73=\"],\n  [\"396\u00f73=\", \"342\u00f78=\"],\n  [\"866\u00f76=\", \"308\u00f75=\"],\n  [\"302\u00f73=\", \"140\u00f76=\"],\n  [\"859\u00f75=\", \"990\u00f79=\"],\n  [\"838\u00f79=\", \"784\u00f75=\"],\n  [\"572\u00f76=\", \"662\u00f79=\"],\n  [\"652\u00f78=\", \"947\u00f78=\"],\n  [\"368\u00f78=\", \"214\u00f78=\"],\n  [\"833\u00f76=\", \"236\u00f78=\"],\n  [\"793\u00f72=\", \"508\u00f72=\"],\n  [\"150\u00f72=\", \"899\u00f74=\"],\n  [\"116\u00f77=\", \"482\u00f76=\"],\n  [\"916\u00f79=\", \"437\u00f77=\"],\n  [\"479\u00f76=\", \"291\u00f73=\"],\n  [\"651\u00f72=\", \"299\u00f74=\"],\n  [\"530\u00f79=\", \"667\u00f74=\"],\n  [\"130\u00f74=\", \"292\u00f76=\"],\n  [\"889\u00f78=\", \"203\u00f72=\"],\n  [\"331\u00f73=\", \"882\u00f78=\"],\n  [\"523\u00f75=\", \"731\u00f77=\"],\n  [\"519\u00f73=\", \"589\u00f77=\"],\n  [\"550\u00f74=\", \"201\u00f77=\"],\n  [\"290\u00f75=\", \"957\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 division problems to the new values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-01-27 Tuesday\", \"2026-01-28 Wednesday\"),\n    @(\"441\u00f72=\", \"143\u00f78=\"),\n    @(\"145\u00f72=\", \"808\u00f73=\"),\n    @(\"396\u00f73=\", \"342\u00f78=\"),\n    @(\"866\u00f76=\", \"308\u00f75=\"),\n    @(\"302\u00f73=\", \"140\u00f76=\"),\n    @(\"859\u00f75=\", \"990\u00f79=\"),\n    @(\"838\u00f79=\", \"784\u00f75=\"),\n    @(\"572\u00f76=\", \"662\u00f79=\"),\n    @(\"652\u00f78=\", \"947\u00f78=\"),\n    @(\"368\u00f78=\", \"214\u00f78=\"),\n    @(\"833\u00f76=\", \"236\u00f78=\"),\n    @(\"793\u00f72=\", \"508\u00f72=\"),\n    @(\"150\u00f72=\", \"899\u00f74=\"),\n    @(\"116\u00f77=\", \"482\u00f76=\"),\n    @(\"916\u00f79=\", \"437\u00f77=\"),\n    @(\"479\u00f76=\", \"291\u00f73=\"),\n    @(\"651\u00f72=\", \"299\u00f74=\"),\n    @(\"530\u00f79=\", \"667\u00f74=\"),\n    @(\"130\u00f74=\", \"292\u00f76=\"),\n    @(\"889\u00f78=\", \"203\u00f72=\"),\n    @(\"331\u00f73=\", \"882\u00f78=\"),\n    @(\"523\u00f75=\", \"731\u00f77=\"),\n    @(\"519\u00f73=\", \"589\u00f77=\"),\n    @(\"550\u00f74=\", \"201\u00f77=\"),\n    @(\"290\u00f75=\", \"957\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
